# Generate Report for Handback
#
# For each localized sheet (zh-cn, de-de) the handoff report gains two new
# populated columns - "Latest Target File" (E) and "Latest Handback File"
# (F) - mirroring the existing "Source File Name" (A) / "Latest Handoff
# File" (C) hyperlinks, the "Latest Handback DateTime" (G) placeholder is
# replaced with the real handback timestamp, and the "Status" column flips
# from "Ready for handoff" to "Handed back: in sync with en-US" everywhere
# that text appears (Overview sheet included, since it mirrors the same
# status text).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# Cornflowerblue (FF6495ED), matching the workbook's existing "HyperLink"
# cell style, expressed as an OLE/VB color value (0x00BBGGRR).
$hyperlinkColor = 15570276

function Get-HyperlinkUrl($ws, $cellRef) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $cellRef) {
            return $h.Address
        }
    }
    return $null
}

function Style-AsHyperlink($ws, $addr) {
    $ws.Range($addr).Font.Underline = 2
    $ws.Range($addr).Font.Color = $hyperlinkColor
}

function Add-MirrorHyperlink($ws, $srcAddr, $dstAddr, $displayText) {
    $url = Get-HyperlinkUrl $ws $srcAddr
    $ws.Hyperlinks.Add($ws.Range($dstAddr), $url, "", "", $displayText) | Out-Null
    Style-AsHyperlink $ws $dstAddr
}

# ---------------------------------------------------------------------
# Overview sheet: just the status-text rename (columns B & C, rows 2-3).
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# Per-language detail sheets.
# ---------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; HandbackRow2 = "2016-03-08 07:06:11"; HandbackRow3 = "2016-03-08 07:06:11" },
    @{ Name = "de-de"; HandbackRow2 = "2016-03-08 07:06:29"; HandbackRow3 = "2016-03-08 07:06:29" }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Status: "Ready for handoff" -> "Handed back: in sync with en-US"
    $ws.Range("B2").Value = $newStatus
    $ws.Range("B3").Value = $newStatus

    # "Latest Target File" (E) mirrors "Source File Name" (A) hyperlink.
    $aText2 = $ws.Range("A2").Text
    Add-MirrorHyperlink $ws '$A$2' "E2" $aText2

    $aText3 = $ws.Range("A3").Text
    Add-MirrorHyperlink $ws '$A$3' "E3" $aText3

    # "Latest Handback File" (F) mirrors "Latest Handoff File" (C) hyperlink.
    $cText2 = $ws.Range("C2").Text
    Add-MirrorHyperlink $ws '$C$2' "F2" $cText2

    $cText3 = $ws.Range("C3").Text
    Add-MirrorHyperlink $ws '$C$3' "F3" $cText3

    # "Latest Handback DateTime" (G): placeholder -> real handback time.
    $ws.Range("G2").Value = $lang.HandbackRow2
    $ws.Range("G3").Value = $lang.HandbackRow3
}

Write-Output "Generated handback report."
